$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.078.08'
$ws.Range("E2").Value = '  -0.75%  '
$ws.Range("D3").Value = '2.610.26'
$ws.Range("E3").Value = '  -2.37%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.47'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.16'
$ws.Range("E6").Value = '  +0.48%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("D9").Value = '2.610.95'
$ws.Range("E9").Value = '  -2.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.49'
$ws.Range("E11").Value = '  -3.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.371'
$ws.Range("E12").Value = '  +3.94%  '
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.05'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("D15").Value = '3.079.43'
$ws.Range("E15").Value = '  -2.20%  '
$ws.Range("D16").Value = '62.941.82'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("D18").Value = '2.604.58'
$ws.Range("E18").Value = '  -3.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.45'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.49'
$ws.Range("E20").Value = '  +1.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '340.87'
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.84'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.74'
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.04'
$ws.Range("E25").Value = '  -2.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.68'
$ws.Range("E26").Value = '  -0.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.58'
$ws.Range("E27").Value = '  +2.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.97'
$ws.Range("E28").Value = '  +5.06%  '
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '542.02'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  +0.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.78'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.02'
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("D34").Value = '0.0₃0836'
$ws.Range("E34").Value = '  +2.60%  '
$ws.Range("E35").Value = '  -5.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.18'
$ws.Range("E36").Value = '  +0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '168.30'
$ws.Range("E37").Value = '  -2.67%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  +4.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.87'
$ws.Range("E41").Value = '  -1.80%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '164.82'
$ws.Range("E43").Value = '  -5.56%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.62'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.74'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '21.69'
$ws.Range("E46").Value = '  -2.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0560'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("E48").Value = '  -2.12%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0241'
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("E50").Value = '  -1.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.91'
$ws.Range("E51").Value = '  +9.94%  '
